# Update figures list in ms
#
# This script reproduces, via Word COM-interop calls, the OOXML diff that:
#  1. Removes the (hidden) "_GoBack" bookmark from the "Earliest and latest
#     years" bullet.
#  2. Merges the three runs "Only " / "from " / "rodent species" into a
#     single run "Only from rodent species".
#  3. Rewrites "3 panel with yearly temp, yearly mass combined, and mrt
#     combined" into "Yearly temp, yearly mass, mrt combined, and r
#     distribution per site" (re-using the existing "mrt" spell-check
#     exception run) and splits the new leading text across three runs.
#  4. Adds two new bullets ("Mass change over time compared to temp change
#     over time by species" and "Table/plot of ARIMA model p-values") to
#     the figures list.
#  5. Moves the "_GoBack" bookmark to the very last (empty) paragraph of
#     the document.

$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1. Drop the "_GoBack" bookmark that currently sits on the "Earliest and
#    latest years" bullet. (It gets re-created at the end of the document
#    in step 5.)
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. "Only " + "from " + "rodent species" -> "Only from rodent species"
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex("*rodent species*")
$p = $d.Paragraphs($idx)
$full = $d.Range($p.Range.Start, $p.Range.End - 1)
$full.Text = ""
$full.InsertAfter("Only from rodent species")

# ---------------------------------------------------------------------
# 3. Rewrite the "3 panel with yearly temp, ..." bullet.
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex("*3 panel with yearly temp*")
$p = $d.Paragraphs($idx)
$pStart = $p.Range.Start
$oldLead = "3 panel with yearly temp, yearly mass combined, and "

# Insert the three new leading runs, each at the true start of the
# paragraph (back to front) so they come out in reading order as three
# distinct <w:r> elements ahead of the untouched "mrt" run.
$r = $d.Range($pStart, $p.Range.End - 1)
$r.InsertBefore(" ")

$r = $d.Range($pStart, $p.Range.End - 1)
$r.InsertBefore(" temp, yearly mass,")

$r = $d.Range($pStart, $p.Range.End - 1)
$r.InsertBefore("Yearly")

# Remove the now-redundant original leading run's text entirely.
$full = $p.Range.Text
$oldLeadPos = $pStart + $full.IndexOf($oldLead)
$rOld = $d.Range($oldLeadPos, $oldLeadPos + $oldLead.Length)
$rOld.Text = ""

# Append the new trailing clause as its own run.
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.InsertAfter(", and r distribution per site")

# ---------------------------------------------------------------------
# 4. Add the two new figure bullets after the paragraph just edited.
# ---------------------------------------------------------------------
$idx = Find-ParagraphIndex("*Yearly*temp, yearly mass,*mrt combined*")
$pEmpty = $d.Paragraphs($idx + 1)
$pEmpty.Range.InsertAfter("Mass change over time compared to temp change over time by species")

$pEmpty.Range.InsertParagraphAfter()
$idxNew = $idx + 2
$pNew = $d.Paragraphs($idxNew)
$pNew.Range.InsertAfter("Table/plot of ARIMA model p-values")

# ---------------------------------------------------------------------
# 5. Re-create the "_GoBack" bookmark on the last paragraph of the doc.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
